# Adds 21 new sentence/word rows (rows 222-242) to the 'sentences' sheet,
# corresponding to issue #22 new vocabulary entries dated 2020-12-08.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(222, 1).Value = "현미 수술 분야의 개척자"
$ws.Cells.Item(222, 2).Value = "a pioneer in the field of microsurgery"
$ws.Cells.Item(222, 3).Value = "분야"
$ws.Cells.Item(222, 4).NumberFormat = "@"
$ws.Cells.Item(222, 4).Value = "2020-12-08"
$ws.Cells.Item(222, 4).ClearFormats()

$ws.Cells.Item(223, 1).Value = "농업 분야 내부적으로 불만이 있다."
$ws.Cells.Item(223, 2).Value = "There is discontent within the farming industry."
$ws.Cells.Item(223, 3).Value = "분야"
$ws.Cells.Item(223, 4).NumberFormat = "@"
$ws.Cells.Item(223, 4).Value = "2020-12-08"
$ws.Cells.Item(223, 4).ClearFormats()

$ws.Cells.Item(224, 1).Value = "회사가 보험 판매 분야로 새로 진출했다."
$ws.Cells.Item(224, 2).Value = "The company branched out into selling insurance."
$ws.Cells.Item(224, 3).Value = "분야"
$ws.Cells.Item(224, 4).NumberFormat = "@"
$ws.Cells.Item(224, 4).Value = "2020-12-08"
$ws.Cells.Item(224, 4).ClearFormats()

$ws.Cells.Item(225, 1).Value = "아주 수월하게 결승전에 진출하다."
$ws.Cells.Item(225, 2).Value = "coast to the final"
$ws.Cells.Item(225, 3).Value = "진출하다"
$ws.Cells.Item(225, 4).NumberFormat = "@"
$ws.Cells.Item(225, 4).Value = "2020-12-08"
$ws.Cells.Item(225, 4).ClearFormats()

$ws.Cells.Item(226, 1).Value = "멕시코의 수출시장에 진출 하다."
$ws.Cells.Item(226, 2).Value = "make inroads in the export market in Mexico"
$ws.Cells.Item(226, 3).Value = "진출하다"
$ws.Cells.Item(226, 4).NumberFormat = "@"
$ws.Cells.Item(226, 4).Value = "2020-12-08"
$ws.Cells.Item(226, 4).ClearFormats()

$ws.Cells.Item(227, 1).Value = "노력하여 신용과 책임있는 지위에 진출하다."
$ws.Cells.Item(227, 2).Value = "work one's way up to a position of trust and responsibility"
$ws.Cells.Item(227, 3).Value = "진출하다"
$ws.Cells.Item(227, 4).NumberFormat = "@"
$ws.Cells.Item(227, 4).Value = "2020-12-08"
$ws.Cells.Item(227, 4).ClearFormats()

$ws.Cells.Item(228, 1).Value = "구두쇠는 언제나 가난하다. "
$ws.Cells.Item(228, 2).Value = "The miser is always poor."
$ws.Cells.Item(228, 3).Value = "가난하다"
$ws.Cells.Item(228, 4).NumberFormat = "@"
$ws.Cells.Item(228, 4).Value = "2020-12-08"
$ws.Cells.Item(228, 4).ClearFormats()

$ws.Cells.Item(229, 1).Value = "돈이 없는 사람은 가난하다. 돈 밖에 없는 사람은 더 가난하다. "
$ws.Cells.Item(229, 2).Value = "욕망하는 자는 늘 가난하다. `n명언/속담`nThe fellow who has no money is poor; the fellow who has nothing but money is poorer still. "
$ws.Cells.Item(229, 3).Value = "가난하다"
$ws.Cells.Item(229, 4).NumberFormat = "@"
$ws.Cells.Item(229, 4).Value = "2020-12-08"
$ws.Cells.Item(229, 4).ClearFormats()

$ws.Cells.Item(230, 1).Value = "The Smiths do not have a pot to piss in."
$ws.Cells.Item(230, 2).Value = "브라운씨 일가는 아주 가난하다."
$ws.Cells.Item(230, 3).Value = "가난하다"
$ws.Cells.Item(230, 4).NumberFormat = "@"
$ws.Cells.Item(230, 4).Value = "2020-12-08"
$ws.Cells.Item(230, 4).ClearFormats()

$ws.Cells.Item(231, 1).Value = "우리는 신혼여행을 베니스로 갔다."
$ws.Cells.Item(231, 2).Value = "We went to Venice for our honeymoon."
$ws.Cells.Item(231, 3).Value = "신혼여행"
$ws.Cells.Item(231, 4).NumberFormat = "@"
$ws.Cells.Item(231, 4).Value = "2020-12-08"
$ws.Cells.Item(231, 4).ClearFormats()

$ws.Cells.Item(232, 1).Value = "그들은 신혼여행지로 파리를 택했다."
$ws.Cells.Item(232, 2).Value = "They’ve fixed on Paris for their honeymoon."
$ws.Cells.Item(232, 3).Value = "신혼여행"
$ws.Cells.Item(232, 4).NumberFormat = "@"
$ws.Cells.Item(232, 4).Value = "2020-12-08"
$ws.Cells.Item(232, 4).ClearFormats()

$ws.Cells.Item(233, 1).Value = "신혼여행은 어디로 갈 건지 결정했니?"
$ws.Cells.Item(233, 2).Value = "Have you made up your minds where to go for your honeymoon?"
$ws.Cells.Item(233, 3).Value = "신혼여행"
$ws.Cells.Item(233, 4).NumberFormat = "@"
$ws.Cells.Item(233, 4).Value = "2020-12-08"
$ws.Cells.Item(233, 4).ClearFormats()

$ws.Cells.Item(234, 1).Value = "어떤 경우에도 문을 열지 마."
$ws.Cells.Item(234, 2).Value = "Don’t open the door, in any circumstances."
$ws.Cells.Item(234, 3).Value = "경우"
$ws.Cells.Item(234, 4).NumberFormat = "@"
$ws.Cells.Item(234, 4).Value = "2020-12-08"
$ws.Cells.Item(234, 4).ClearFormats()

$ws.Cells.Item(235, 1).Value = "12시간 이상 지체될 경우에는 여비를 전액 환불해 드립니다."
$ws.Cells.Item(235, 2).Value = "If there is a delay of 12 hours or more, you will receive a full refund of the price of your trip."
$ws.Cells.Item(235, 3).Value = "경우"
$ws.Cells.Item(235, 4).NumberFormat = "@"
$ws.Cells.Item(235, 4).Value = "2020-12-08"
$ws.Cells.Item(235, 4).ClearFormats()

$ws.Cells.Item(236, 1).Value = "이 경우는 분명히 법령 2001의 범위 내에 든다."
$ws.Cells.Item(236, 2).Value = "This case falls clearly within the ambit of the 2001 act."
$ws.Cells.Item(236, 3).Value = "경우"
$ws.Cells.Item(236, 4).NumberFormat = "@"
$ws.Cells.Item(236, 4).Value = "2020-12-08"
$ws.Cells.Item(236, 4).ClearFormats()

$ws.Cells.Item(237, 1).Value = "아이들은 수영장에서 서로를 물 속에 밀어 넣으며 장난을 치고 있었다."
$ws.Cells.Item(237, 2).Value = "The kids were ducking each other in the pool."
$ws.Cells.Item(237, 3).Value = "장난"
$ws.Cells.Item(237, 4).NumberFormat = "@"
$ws.Cells.Item(237, 4).Value = "2020-12-08"
$ws.Cells.Item(237, 4).ClearFormats()

$ws.Cells.Item(238, 1).Value = "음식 가지고 장난 그만 해!"
$ws.Cells.Item(238, 2).Value = "Stop playing with your food!"
$ws.Cells.Item(238, 3).Value = "장난"
$ws.Cells.Item(238, 4).NumberFormat = "@"
$ws.Cells.Item(238, 4).Value = "2020-12-08"
$ws.Cells.Item(238, 4).ClearFormats()

$ws.Cells.Item(239, 1).Value = "그것은 모두 악의 없는 장난이었다."
$ws.Cells.Item(239, 2).Value = "It was all innocent fun."
$ws.Cells.Item(239, 3).Value = "장난"
$ws.Cells.Item(239, 4).NumberFormat = "@"
$ws.Cells.Item(239, 4).Value = "2020-12-08"
$ws.Cells.Item(239, 4).ClearFormats()

$ws.Cells.Item(240, 1).Value = "이리 와, 이놈의 장난꾸러기 녀석!"
$ws.Cells.Item(240, 2).Value = "Come here, you cheeky little monkey!"
$ws.Cells.Item(240, 3).Value = "장난꾸러기"
$ws.Cells.Item(240, 4).NumberFormat = "@"
$ws.Cells.Item(240, 4).Value = "2020-12-08"
$ws.Cells.Item(240, 4).ClearFormats()

$ws.Cells.Item(241, 1).Value = "장난꾸러기 학생들이 학교 건물에 상당한 손상을 입혔다."
$ws.Cells.Item(241, 2).Value = "Student pranksters have done considerable damage to the school buildings."
$ws.Cells.Item(241, 3).Value = "장난꾸러기"
$ws.Cells.Item(241, 4).NumberFormat = "@"
$ws.Cells.Item(241, 4).Value = "2020-12-08"
$ws.Cells.Item(241, 4).ClearFormats()

$ws.Cells.Item(242, 1).Value = "내 조카는 장난꾸러기 소년이다."
$ws.Cells.Item(242, 2).Value = "My nephew is a naughty boy."
$ws.Cells.Item(242, 3).Value = "장난꾸러기"
$ws.Cells.Item(242, 4).NumberFormat = "@"
$ws.Cells.Item(242, 4).Value = "2020-12-08"
$ws.Cells.Item(242, 4).ClearFormats()
